$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Text fixes: "Saving" -> "Savings" (4 occurrences) ---

# Shape 2 "Freeform 5": "  [Saving User] ..." -> "  [Savings User] ..."
$sh2 = $s.Shapes.Item(2)
$run2 = $sh2.TextFrame.TextRange.Characters(1, 15)
$run2.Text = "  [Savings User]"

# Shape 8 "Freeform 37": "[Saving User] ..." -> "[Savings User] ..."
$sh8 = $s.Shapes.Item(8)
$run8 = $sh8.TextFrame.TextRange.Characters(1, 13)
$run8.Text = "[Savings User]"

# Shape 11 "TextBox 21": "[Saving Group Admin] ..." -> "[Savings Group Admin] ..."
$sh11 = $s.Shapes.Item(11)
$run11 = $sh11.TextFrame.TextRange.Characters(1, 20)
$run11.Text = "[Savings Group Admin]"

# Shape 12 "TextBox 24": "[Saving Group Admin] ..." -> "[Savings Group Admin] ..."
$sh12 = $s.Shapes.Item(12)
$run12 = $sh12.TextFrame.TextRange.Characters(1, 20)
$run12.Text = "[Savings Group Admin]"

# --- Position nudges for the three small connector ovals ---

# Shape 15 "Oval 2": (183356, 1136350) -> (172723, 1168249) EMU
$sh15 = $s.Shapes.Item(15)
$sh15.Left = 13.600275590551181
$sh15.Top = 91.98814960629922

# Shape 16 "Oval 27": (1720048, 2499015) -> (1709415, 2477749) EMU
$sh16 = $s.Shapes.Item(16)
$sh16.Left = 134.59964566929133
$sh16.Top = 195.09838582677165

# Shape 17 "Oval 28": (4387346, 2031518) -> (4387346, 2052784) EMU
$sh17 = $s.Shapes.Item(17)
$sh17.Left = 345.46035433070864
$sh17.Top = 161.6365748031496
